$wb = $excel.ActiveWorkbook

# Each worksheet (FE_LFT_#1, FE_LFT_#2, FE_PLT_#1, FE_PLT_#2) gets one new
# daily-log row appended (row 68), mirroring the formatting of the last
# existing row (row 67) and then filling in the new record's values.

$rows = @(
    @{
        Sheet = "FE_LFT_#1"
        A = [double]"45854.49613425926"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x3C"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 316
        I = 15
    },
    @{
        Sheet = "FE_LFT_#2"
        A = [double]"45854.49613425926"
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x4C"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 332
        I = 14
    },
    @{
        Sheet = "FE_PLT_#1"
        A = [double]"45854.49613425926"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x63"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 99
        I = 3
    },
    @{
        Sheet = "FE_PLT_#2"
        A = [double]"45854.49613425926"
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x63"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 99
        I = 3
    }
)

foreach ($rec in $rows) {
    $ws = $wb.Worksheets.Item($rec.Sheet)

    # Seed row 68 with the formatting (number formats, styles) of row 67,
    # then overwrite with the new record's actual values.
    $ws.Range("A67:I67").Copy($ws.Range("A68:I68"))

    $ws.Range("A68").Value = $rec.A
    $ws.Range("B68").Value = $rec.B
    $ws.Range("C68").Value = $rec.C
    $ws.Range("D68").Value = $rec.D
    $ws.Range("E68").Value = $rec.E
    $ws.Range("F68").Value = $rec.F
    $ws.Range("G68").Value = $rec.G
    $ws.Range("H68").Value = $rec.H
    $ws.Range("I68").Value = $rec.I
}
